$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D4 gets new text (a new shared string, "“dtyui”").
$ws.Range("D4").Value = "“dtyui”"

# The active selection moves from C4 to D4.
$ws.Range("D4").Select()
